# data change / new graphs
#
# 1. Rename the sheet from "Sheet1" to "precincts".
# 2. Remove the three districtr.org COI hyperlinks (and their display text)
#    that lived in B193:B195 -- they become empty (but still styled) cells.
#    Excel prunes the now-unused shared strings automatically on save.
# 3. Update the view: the user had scrolled/selected down near the bottom of
#    the sheet (selection moved from C66 to C144, and the frozen top pane's
#    scroll position moved down to row 122).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "precincts"

# 2. Drop the three hyperlinks + their cell text
$ws.Range("B193").Hyperlinks.Delete() | Out-Null
$ws.Range("B194").Hyperlinks.Delete() | Out-Null
$ws.Range("B195").Hyperlinks.Delete() | Out-Null
$ws.Range("B193:B195").ClearContents() | Out-Null

# 3. Move the selection / scroll position
$ws.Range("C144").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 122
$win.ScrollColumn = 1
